# Append a new paragraph "Mailo parte3" after "Chao mundo parte2".
# "Mailo" is flagged by the spell checker (proofErr spellStart/spellEnd)
# and the text is split into two runs ("Mailo" and " parte3"), matching
# how Word records a freshly-typed, proofed word.

$d = $word.ActiveDocument

# The existing last paragraph is "Chao mundo parte2"; insert a brand new
# paragraph right after it. InsertParagraphAfter() on its Range creates a
# new paragraph that inherits the same paragraph/run formatting
# (w:rPr/w:lang es-ES) as the paragraph it follows.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter() | Out-Null

# Grab the freshly created (now last) paragraph and fill it with the
# exact run/proofErr structure from the target markup.
$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range

$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Mailo</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> parte3</w:t></w:r>' +
  '</w:p>'

$newRange.InsertXML($xmlFrag) | Out-Null
